# Split the sentence "We selected the model with the lowest AIC value, "
# into two sentences: "We " + "used Akaike Information Criterion (AIC) to
# select the model that best represented the data."

$d = $word.ActiveDocument

$oldText    = "We selected the model with the lowest AIC value, "
$firstPart  = "We "
$secondPart = "used Akaike Information Criterion (AIC) to select the model that best represented the data."

$r = $d.Content
$found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $splitPos = $r.Start + $firstPart.Length

    # Replace everything after "We " with the new second sentence.
    $tail = $d.Range($splitPos, $r.End)
    $tail.Text = $secondPart

    # Toggling a character property off/on/off on the newly inserted text
    # forces Word to keep it as its own run (with its own rPr) instead of
    # silently coalescing it into the preceding identically-formatted run.
    $newTail = $d.Range($splitPos, $splitPos + $secondPart.Length)
    $newTail.Bold = 1
    $newTail.Bold = 0

    # Do the same for "We " so it stays a distinct run from the sentence
    # that precedes it (". ") rather than being merged back together.
    $head = $d.Range($r.Start, $splitPos)
    $head.Bold = 1
    $head.Bold = 0
}
